# Update "想去人数" (number of people wanting to go) values in column F
# across sheets "展览" (sheet1), "演出" (sheet2), and "全部类型" (sheet4)
# to reflect newly scraped counts, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 79
$wsExhibit.Range("F4").Value = 720
$wsExhibit.Range("F7").Value = 2796
$wsExhibit.Range("F8").Value = 1674
$wsExhibit.Range("F9").Value = 1792
$wsExhibit.Range("F11").Value = 283
$wsExhibit.Range("F12").Value = 727
$wsExhibit.Range("F13").Value = 885
$wsExhibit.Range("F15").Value = 368
$wsExhibit.Range("F16").Value = 1112
$wsExhibit.Range("F20").Value = 6513
$wsExhibit.Range("F21").Value = 251
$wsExhibit.Range("F22").Value = 1467
$wsExhibit.Range("F23").Value = 149
$wsExhibit.Range("F26").Value = 306
$wsExhibit.Range("F27").Value = 261
$wsExhibit.Range("F30").Value = 898
$wsExhibit.Range("F34").Value = 470
$wsExhibit.Range("F35").Value = 1387
$wsExhibit.Range("F36").Value = 157
$wsExhibit.Range("F38").Value = 216
$wsExhibit.Range("F39").Value = 8
$wsExhibit.Range("F41").Value = 192
$wsExhibit.Range("F42").Value = 154

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F9").Value = 7

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 79
$wsAll.Range("F4").Value = 720
$wsAll.Range("F10").Value = 2796
$wsAll.Range("F11").Value = 1674
$wsAll.Range("F12").Value = 1792
$wsAll.Range("F14").Value = 283
$wsAll.Range("F15").Value = 727
$wsAll.Range("F17").Value = 885
$wsAll.Range("F19").Value = 368
$wsAll.Range("F20").Value = 1112
$wsAll.Range("F23").Value = 6513
$wsAll.Range("F24").Value = 251
$wsAll.Range("F25").Value = 1467
$wsAll.Range("F27").Value = 149
$wsAll.Range("F30").Value = 306
$wsAll.Range("F31").Value = 261
$wsAll.Range("F34").Value = 898
$wsAll.Range("F38").Value = 470
$wsAll.Range("F39").Value = 1387
$wsAll.Range("F40").Value = 157
$wsAll.Range("F42").Value = 216
$wsAll.Range("F43").Value = 8
$wsAll.Range("F45").Value = 192
$wsAll.Range("F48").Value = 7
$wsAll.Range("F49").Value = 154

